$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric so Excel
# keeps them as literal text (matching the original inlineStr content)
foreach ($addr in @("D5", "D10", "D13", "D15", "D17", "D19", "D22", "D23", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.447.60"
$ws.Range("E2").Value = "  +2.91%  "
$ws.Range("D3").Value = "2.311.48"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "311.67"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("E6").Value = "  +5.86%  "
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +7.30%  "
$ws.Range("D10").Value = "35.84"
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "7.01"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "2.670.18"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "15.02"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "2.306.76"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "0.810"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "43.349.39"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "12.37"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  +3.20%  "
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("D22").Value = "68.12"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "241.39"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "24.65"
$ws.Range("E27").Value = "  +4.72%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "37.15"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.63"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "168.25"
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.29"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.51"
$ws.Range("E34").Value = "  +5.81%  "
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "17.63"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "3.07"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.106"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.116"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "4.38"
$ws.Range("E41").Value = "  +7.87%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "2.32"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "19.40"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0289"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "1.975.15"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "9.91"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "55.51"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.92"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.57"
$ws.Range("E50").Value = "  +7.47%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.538.28"
$ws.Range("E51").Value = "  +1.91%  "
